$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh prepends two new report rows (new week, row 2:3) and
# pushes all previously existing data rows down by two (old row 2 -> row 4,
# ..., old row 12 -> row 14). Insert two blank rows at row 2 to reproduce
# that shift while preserving formatting/styles on the shifted rows.
$ws.Rows("2:3").Insert()

# Newly inserted rows pick up the formatting of the row above (the header),
# which is not what any of the data rows use. Clear that so the new rows
# match the plain (unstyled) look of the other data rows.
$ws.Range("A2:T3").ClearFormats()

# Re-apply the date number format used by the rest of column D to the two
# new date cells.
$ws.Range("D2").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat

# Row 2: Comercializadora del Agro de Limarí, Damasco, Castle Brite, Especial
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44552
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103003
$ws.Range("J2").Value = "Damasco"
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 360
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("Q2").Value = "$/caja 18 kilos"
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1139
$ws.Range("T2").Value = 18

# Row 3: Comercializadora del Agro de Limarí, Damasco, Castle Brite, Primera
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44552
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103003
$ws.Range("J3").Value = "Damasco"
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 280
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 19000
$ws.Range("P3").Value = 18500
$ws.Range("Q3").Value = "$/caja 18 kilos"
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 1028
$ws.Range("T3").Value = 18
